$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 279, shifting existing rows 279:370 down to 280:371
$ws.Rows.Item(279).Insert()

# Populate the newly inserted row 279 with the new weekly price record
$ws.Cells.Item(279, 1).Value = 6
$ws.Cells.Item(279, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(279, 3).Value = "Metropolitana"
$ws.Cells.Item(279, 4).Value = 45215
$ws.Cells.Item(279, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(279, 5).Value = 13
$ws.Cells.Item(279, 6).Value = 100112001
$ws.Cells.Item(279, 7).Value = "Berenjena"
$ws.Cells.Item(279, 8).Value = "Sin especificar"
$ws.Cells.Item(279, 9).Value = "Primera"
$ws.Cells.Item(279, 10).Value = 250
$ws.Cells.Item(279, 11).Value = 8000
$ws.Cells.Item(279, 12).Value = 9000
$ws.Cells.Item(279, 13).Value = 8400
$ws.Cells.Item(279, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(279, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(279, 16).Value = 168
$ws.Cells.Item(279, 17).Value = 50
$ws.Cells.Item(279, 18).Value = "Hortaliza"
